# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 125 on the "Pomelo" sheet,
# shifting all subsequent rows (old 125-145) down by one (new 126-146).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 125 (pushes old row 125 -> 126, etc.)
$ws.Rows("125:125").Insert()

# Populate the new row 125 with the new weekly record
$ws.Range("A125").Value2 = 10
$ws.Range("B125").Value2 = "Vega Modelo de Temuco"
$ws.Range("C125").Value2 = "La Araucanía"
$ws.Range("D125").Value2 = 44474
$ws.Range("E125").Value2 = 9
$ws.Range("F125").Value2 = "Fruta"
$ws.Range("G125").Value2 = 100102
$ws.Range("H125").Value2 = "Cítricos"
$ws.Range("I125").Value2 = 100102006
$ws.Range("J125").Value2 = "Pomelo"
$ws.Range("K125").Value2 = "Start Ruby"
$ws.Range("L125").Value2 = "Especial"
$ws.Range("M125").Value2 = 50
$ws.Range("N125").Value2 = 16000
$ws.Range("O125").Value2 = 16000
$ws.Range("P125").Value2 = 16000
$ws.Range("Q125").Value2 = "$/caja 14 kilos granel"
$ws.Range("R125").Value2 = "Región de O'Higgins"
$ws.Range("S125").Value2 = 1143
$ws.Range("T125").Value2 = 14
